$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) of the last existing row down to the new rows
$ws.Range("A464:D464").Copy($ws.Range("A465:D491"))

# Fill in the new data rows (9 Dec 2021 .. 5 Jan 2022)
$ws.Range("A465").Value = 44539
$ws.Range("B465").Value = 14
$ws.Range("C465").Value = 86
$ws.Range("D465").Value = 356.5210181576983
$ws.Range("A466").Value = 44540
$ws.Range("B466").Value = 4
$ws.Range("C466").Value = 56
$ws.Range("D466").Value = 232.1532211259431
$ws.Range("A467").Value = 44541
$ws.Range("B467").Value = 8
$ws.Range("C467").Value = 63
$ws.Range("D467").Value = 261.172373766686
$ws.Range("A468").Value = 44542
$ws.Range("B468").Value = 18
$ws.Range("C468").Value = 72
$ws.Range("D468").Value = 298.4827128762126
$ws.Range("A469").Value = 44543
$ws.Range("B469").Value = 22
$ws.Range("C469").Value = 85
$ws.Range("D469").Value = 352.3754249233065
$ws.Range("A470").Value = 44544
$ws.Range("B470").Value = 9
$ws.Range("C470").Value = 75
$ws.Range("D470").Value = 310.9194925793881
$ws.Range("A471").Value = 44545
$ws.Range("B471").Value = 2
$ws.Range("C471").Value = 77
$ws.Range("D471").Value = 319.2106790481718
$ws.Range("A472").Value = 44546
$ws.Range("B472").Value = 10
$ws.Range("C472").Value = 73
$ws.Range("D472").Value = 302.6283061106044
$ws.Range("A473").Value = 44547
$ws.Range("B473").Value = 5
$ws.Range("C473").Value = 74
$ws.Range("D473").Value = 306.7738993449963
$ws.Range("A474").Value = 44548
$ws.Range("B474").Value = 11
$ws.Range("C474").Value = 77
$ws.Range("D474").Value = 319.2106790481718
$ws.Range("A475").Value = 44550
$ws.Range("B475").Value = 12
$ws.Range("C475").Value = 71
$ws.Range("D475").Value = 294.3371196418207
$ws.Range("A476").Value = 44551
$ws.Range("B476").Value = 14
$ws.Range("C476").Value = 63
$ws.Range("D476").Value = 261.172373766686
$ws.Range("A477").Value = 44552
$ws.Range("B477").Value = 0
$ws.Range("C477").Value = 54
$ws.Range("D477").Value = 223.8620346571595
$ws.Range("A478").Value = 44553
$ws.Range("B478").Value = 16
$ws.Range("C478").Value = 68
$ws.Range("D478").Value = 281.9003399386452
$ws.Range("A479").Value = 44554
$ws.Range("B479").Value = 1
$ws.Range("C479").Value = 59
$ws.Range("D479").Value = 244.5900008291186
$ws.Range("A480").Value = 44555
$ws.Range("B480").Value = 5
$ws.Range("C480").Value = 59
$ws.Range("D480").Value = 244.5900008291186
$ws.Range("A481").Value = 44556
$ws.Range("B481").Value = 19
$ws.Range("C481").Value = 67
$ws.Range("D481").Value = 277.7547467042534
$ws.Range("A482").Value = 44557
$ws.Range("B482").Value = 26
$ws.Range("C482").Value = 81
$ws.Range("D482").Value = 335.7930519857392
$ws.Range("A483").Value = 44558
$ws.Range("B483").Value = 42
$ws.Range("C483").Value = 109
$ws.Range("D483").Value = 451.8696625487107
$ws.Range("A484").Value = 44559
$ws.Range("B484").Value = 13
$ws.Range("C484").Value = 122
$ws.Range("D484").Value = 505.7623745958047
$ws.Range("A485").Value = 44560
$ws.Range("B485").Value = 14
$ws.Range("C485").Value = 120
$ws.Range("D485").Value = 497.471188127021
$ws.Range("A486").Value = 44561
$ws.Range("B486").Value = 29
$ws.Range("C486").Value = 148
$ws.Range("D486").Value = 613.5477986899925
$ws.Range("A487").Value = 44562
$ws.Range("B487").Value = 16
$ws.Range("C487").Value = 159
$ws.Range("D487").Value = 659.1493242683027
$ws.Range("A488").Value = 44563
$ws.Range("B488").Value = 47
$ws.Range("C488").Value = 187
$ws.Range("D488").Value = 775.2259348312743
$ws.Range("A489").Value = 44564
$ws.Range("B489").Value = 47
$ws.Range("C489").Value = 208
$ws.Range("D489").Value = 862.283392753503
$ws.Range("A490").Value = 44565
$ws.Range("B490").Value = 25
$ws.Range("C490").Value = 191
$ws.Range("D490").Value = 791.8083077688417
$ws.Range("A491").Value = 44566
$ws.Range("B491").Value = 32
$ws.Range("C491").Value = 210
$ws.Range("D491").Value = 870.5745792222866

